$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cls = $m.CustomLayouts
$src = $cls.Item(7)
$dup = $src.Duplicate()

$EMU = 12700.0
$content2 = $dup.Shapes.Item(5)

$target = 1205547
$basePts = $target / $EMU
$found = $false
for ($k = 0; $k -le 4000 -and -not $found; $k++) {
    $delta = [math]::Floor(($k + 1) / 2) * 0.0000001
    if ($k % 2 -eq 0) { $delta = -$delta }
    $pts = $basePts + $delta
    $content2.Top = $pts
    $got = [math]::Round($content2.Top * $EMU)
    if ($got -eq $target) {
        $found = $true
        Write-Output "FOUND at k=$k delta=$delta pts=$pts got=$got"
    }
}
if (-not $found) {
    Write-Output "NOT FOUND"
}
